# Charles_Frank_Cycle13_Presentation.pptx edit
# 1) Update the auto date placeholder text (datetimeFigureOut cache) on the
#    slide master and every slide layout from 11/11/16 -> 11/12/16.
# 2) Move slide 12's lone picture onto slide 11 as a fourth ("Picture 4")
#    image filling the bottom-right quadrant, then delete slide 12.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $isDate = $false
            try {
                if ($sh.PlaceholderFormat.Type -eq 16) { $isDate = $true }
            } catch {
                $isDate = $false
            }
            if ($isDate) {
                if ($sh.TextFrame.TextRange.Text -eq "11/11/16") {
                    $sh.TextFrame.TextRange.Text = "11/12/16"
                }
            }
        }
    }
}

# Slide master date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every custom (slide) layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $cl = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $cl.Shapes
}

# --- Move slide 12's picture onto slide 11, then remove slide 12. ---

$s11 = $p.Slides.Item(11)
$s12 = $p.Slides.Item(12)

# The engine hands out new shape ids from a small per-slide counter (2, 3, 4, ...)
# skipping any id already used on that slide. Slide 11 already has ids 1,4,9,10,
# so the first two free slots are 2 and 3; burn through those with throw-away
# shapes so the picture we paste next lands on id 5 (matching "Picture 4").
$burn1 = $s11.Shapes.AddTextbox(1, 0, 0, 1, 1)
$burn2 = $s11.Shapes.AddTextbox(1, 0, 0, 1, 1)
$burn1.Delete()
$burn2.Delete()

$srcPic = $s12.Shapes.Item(1)
$srcPic.Copy()
$newPic = $s11.Shapes.Paste().Item(1)
$newPic.Name = "Picture 4"
$newPic.Left = 479.04331970230976
$newPic.Top = 260.8162994389098
$newPic.Width = 489.75
$newPic.Height = 279.18370078740156

$s12.Delete()
